# Add new header columns (L:T) + new data row values (L2:T2), and
# apply center-alignment formatting to the whole header row (A1:T1),
# matching the "Add files via upload" revision of SMDB.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (columns L..T on row 1) ------------------------------
$ws.Range("L1").Value = "MD"
$ws.Range("M1").Value = "MP"
$ws.Range("N1").Value = "MW"
$ws.Range("O1").Value = "MB"
$ws.Range("P1").Value = "MR"
$ws.Range("Q1").Value = "MC"
$ws.Range("R1").Value = "MT"
$ws.Range("S1").Value = "Ml"
$ws.Range("T1").Value = "MA"

# Give the new header cells the same formatting as the existing header
# cells (font2, no fill) before the across-the-board center alignment below.
$ws.Range("K1").Copy()
$ws.Range("L1:T1").PasteSpecial(-4122)

# --- New data cells (columns L..T on row 2) --------------------------------
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 3

# --- Center-align the full header row, A1:T1 --------------------------------
$xlCenter = -4108
$ws.Range("A1:T1").HorizontalAlignment = $xlCenter

# --- Restore the selection left behind by the edit --------------------------
$ws.Range("I17").Select() | Out-Null
